$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 78
$ws.Range("H2").Value = 84
$ws.Range("E3").Value = 45
$ws.Range("F3").Value = 33
$ws.Range("H3").Value = 34
$ws.Range("E4").Value = 49
$ws.Range("F4").Value = 30
$ws.Range("H4").Value = 42
$ws.Range("F5").Value = 108
$ws.Range("H5").Value = 119
$ws.Range("E6").Value = 49
$ws.Range("F6").Value = 35
$ws.Range("H6").Value = 45
$ws.Range("F7").Value = 28
$ws.Range("H7").Value = 32
$ws.Range("F8").Value = 6
$ws.Range("H8").Value = 9
$ws.Range("E10").Value = 658
$ws.Range("F10").Value = 368
$ws.Range("H10").Value = 464
$ws.Range("E11").Value = 432
$ws.Range("F11").Value = 245
$ws.Range("H11").Value = 309
$ws.Range("E12").Value = 660
$ws.Range("F12").Value = 403
$ws.Range("H12").Value = 489
$ws.Range("F13").Value = 91
$ws.Range("H13").Value = 125
$ws.Range("E14").Value = 136
$ws.Range("F14").Value = 80
$ws.Range("H14").Value = 114
$ws.Range("E15").Value = 192
$ws.Range("F15").Value = 91
$ws.Range("H15").Value = 141
$ws.Range("E16").Value = 227
$ws.Range("F16").Value = 132
$ws.Range("H16").Value = 180
$ws.Range("F17").Value = 66
$ws.Range("H17").Value = 90
$ws.Range("F19").Value = 9
$ws.Range("H19").Value = 12
$ws.Range("F20").Value = 40
$ws.Range("H20").Value = 77
$ws.Range("E22").Value = 188
$ws.Range("F22").Value = 105
$ws.Range("H22").Value = 147
$ws.Range("E23").Value = 218
$ws.Range("F23").Value = 112
$ws.Range("H23").Value = 163
$ws.Range("E24").Value = 249
$ws.Range("F24").Value = 149
$ws.Range("H24").Value = 179
$ws.Range("E25").Value = 309
$ws.Range("F25").Value = 173
$ws.Range("H25").Value = 233
$ws.Range("E26").Value = 181
$ws.Range("F26").Value = 115
$ws.Range("H26").Value = 140
$ws.Range("F27").Value = 200
$ws.Range("H27").Value = 281
$ws.Range("F28").Value = 108
$ws.Range("H28").Value = 160
$ws.Range("E29").Value = 188
$ws.Range("F29").Value = 114
$ws.Range("H29").Value = 155
$ws.Range("E30").Value = 244
$ws.Range("F30").Value = 154
$ws.Range("H30").Value = 206
$ws.Range("F32").Value = 132
$ws.Range("H32").Value = 170
$ws.Range("E33").Value = 320
$ws.Range("F33").Value = 177
$ws.Range("H33").Value = 267
$ws.Range("F34").Value = 171
$ws.Range("H34").Value = 209
$ws.Range("E35").Value = 173
$ws.Range("F35").Value = 123
$ws.Range("H35").Value = 150
$ws.Range("E36").Value = 88
$ws.Range("F36").Value = 57
$ws.Range("H36").Value = 67
$ws.Range("E37").Value = 185
$ws.Range("F37").Value = 105
$ws.Range("H37").Value = 142
$ws.Range("F38").Value = 62
$ws.Range("H38").Value = 78
$ws.Range("E39").Value = 193
$ws.Range("F39").Value = 102
$ws.Range("H39").Value = 153
$ws.Range("E40").Value = 291
$ws.Range("F40").Value = 149
$ws.Range("H40").Value = 229
$ws.Range("E41").Value = 426
$ws.Range("F41").Value = 213
$ws.Range("H41").Value = 305
$ws.Range("E42").Value = 434
$ws.Range("F42").Value = 252
$ws.Range("H42").Value = 313
$ws.Range("E43").Value = 138
$ws.Range("F43").Value = 76
$ws.Range("H43").Value = 103
$ws.Range("F44").Value = 184
$ws.Range("H44").Value = 252
$ws.Range("E45").Value = 172
$ws.Range("F45").Value = 97
$ws.Range("H45").Value = 136
$ws.Range("E46").Value = 375
$ws.Range("F46").Value = 219
$ws.Range("H46").Value = 282
$ws.Range("E47").Value = 520
$ws.Range("F47").Value = 293
$ws.Range("H47").Value = 385
$ws.Range("E48").Value = 253
$ws.Range("F48").Value = 121
$ws.Range("H48").Value = 165
$ws.Range("E49").Value = 329
$ws.Range("F49").Value = 166
$ws.Range("H49").Value = 253
$ws.Range("E50").Value = 268
$ws.Range("F50").Value = 145
$ws.Range("H50").Value = 216
$ws.Range("E51").Value = 259
$ws.Range("F51").Value = 129
$ws.Range("H51").Value = 203
